$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style s="4" with border) from row 42 down into rows 43-48
$ws.Range("A42:G42").Copy()
$ws.Range("A43:G48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill cell values in the exact order the strings were first authored,
# so the shared-strings table is built up in the same sequence as the source edit.
$ws.Range("A48").Value = "Tutorial text"
$ws.Range("B43").Value = "tutorial wave"
$ws.Range("B44").Value = "tutorial skills"
$ws.Range("B45").Value = "tutorial currency"
$ws.Range("B46").Value = "tutorial chest"
$ws.Range("B47").Value = "tutorial potion"
$ws.Range("B48").Value = "tutorial boss portal"
$ws.Range("C43").Value = "Information on the current wave is located here: timer, name, etc."
$ws.Range("C44").Value = "Your skills. Some can be used more often than others."
$ws.Range("C45").Value = "Your currency. Use it to open chests!"
$ws.Range("C46").Value = "A bonus chest. Open it to grow more powerful!"
$ws.Range("C47").Value = "One potion a day keeps the doctor away."
$ws.Range("E43").Value = "La información sobre la ola actual se encuentra aquí: temporizador, nombre, etc."
$ws.Range("G43").Value = "当前波形的信息位于此处：计时器、名称等。"
$ws.Range("F43").Value = "現在のウェーブに関する情報はここにあります: タイマー、名前など。"
$ws.Range("G44").Value = "你的技能。有些可以比其他更频繁地使用。"
$ws.Range("F44").Value = "あなたの技術。他のものより頻繁に使用できるものもあります。"
$ws.Range("E44").Value = "Tus habilidades. Algunos se pueden usar con más frecuencia que otros."
$ws.Range("E45").Value = "Tu moneda. ¡Úsalo para abrir cofres!"
$ws.Range("F45").Value = "あなたの通貨。それを使って宝箱を開けよう！"
$ws.Range("G45").Value = "你的货币。用它来打开箱子！"
$ws.Range("E46").Value = "Un cofre de bonificación. ¡Ábrelo para volverte más poderoso!"
$ws.Range("G46").Value = "一个奖金箱子。打开它变得更强大！"
$ws.Range("F46").Value = "ボーナスチェスト。それを開いて、より強力に成長してください！"
$ws.Range("G47").Value = "一天一剂，医生远离我。"
$ws.Range("F47").Value = "1日1錠で医者いらず。"
$ws.Range("E47").Value = "Una poción al día mantiene alejado al médico."
$ws.Range("G48").Value = "生成关卡 Boss 的传送门。底部的数字表示打开它并召唤 Boss 所需的货币。"
$ws.Range("F48").Value = "レベルのボスを生成するポータル。下部の数字は、それを開いてボスを呼び出すために必要な通貨を示しています。"
$ws.Range("E48").Value = "Un portal que genera el jefe del nivel. El número en la parte inferior indica la moneda requerida para abrirlo y llamar al jefe."
$ws.Range("D43").Value = "De l'information sur la vague en cours est affichée ici: temps, nom, etc."
$ws.Range("D44").Value = "Vos aptitudes. Certaines peuvent être utilisées plus fréquemment que d'autres."
$ws.Range("D45").Value = "Votre monnaie. Utilisez-la pour ouvrir des coffres!"
$ws.Range("D46").Value = "Un coffre bonus. Ouvrez-le pour devenir plus puissant!"
$ws.Range("D47").Value = "Une potion par jour garde le médecin au loin!"
$ws.Range("D48").Value = "Un portail qui fait apparaître le boss du niveau. Le nombre au bas indique la monnaie nécessaire pour ouvrir le portail et invoquer le boss."
$ws.Range("C48").Value = "A portal that spawns the level's boss.\nThe number at the bottom indicates the required currency to open it and call forth the boss."

# Remaining cells reuse the already-created "Tutorial text" shared string
$ws.Range("A43").Value = "Tutorial text"
$ws.Range("A44").Value = "Tutorial text"
$ws.Range("A45").Value = "Tutorial text"
$ws.Range("A46").Value = "Tutorial text"
$ws.Range("A47").Value = "Tutorial text"

# Row heights (wrap-text auto-fit heights captured from the authored session)
$ws.Range("A43:G43").RowHeight = 100.8
$ws.Range("A44:G44").RowHeight = 100.8
$ws.Range("A45:G45").RowHeight = 72
$ws.Range("A46:G46").RowHeight = 86.4
$ws.Range("A47:G47").RowHeight = 57.6
$ws.Range("A48:G48").RowHeight = 158.4

# C44 keeps the un-bordered default style, matching the source session
$ws.Range("C44").Borders.LineStyle = -4142

# Restore selection state to reflect final editing position
$ws.Range("C48").Select() | Out-Null
